$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44189
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 23500
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 23750
$ws.Range("Q2").Value = "`$/caja 18 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1319
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44189
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 21500
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21750
$ws.Range("Q3").Value = "`$/caja 18 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1208
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44566
$ws.Range("K4").Value = "Modesto"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 23000
$ws.Range("O4").Value = 24000
$ws.Range("P4").Value = 23500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1306

# Row 5
$ws.Range("D5").Value = 44566
$ws.Range("K5").Value = "Modesto"
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1194

# Row 7
$ws.Range("D7").Value = 44546
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 22500
$ws.Range("O7").Value = 23000
$ws.Range("P7").Value = 22750
$ws.Range("Q7").Value = "`$/caja 18 kilos"
$ws.Range("S7").Value = 1264
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44546
$ws.Range("K8").Value = "Castle Brite"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 20500
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20750
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1153

# Row 9
$ws.Range("D9").Value = 44161
$ws.Range("K9").Value = "Dina"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20500
$ws.Range("P9").Value = 20250
$ws.Range("Q9").Value = "`$/caja 15 kilos"
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 1350
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("D10").Value = 44161
$ws.Range("K10").Value = "Dina"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18500
$ws.Range("P10").Value = 18250
$ws.Range("Q10").Value = "`$/caja 15 kilos"
$ws.Range("S10").Value = 1217
$ws.Range("T10").Value = 15

# Row 11
$ws.Range("D11").Value = 44160
$ws.Range("M11").Value = 240
$ws.Range("Q11").Value = "`$/caja 15 kilos"
$ws.Range("S11").Value = 1383
$ws.Range("T11").Value = 15

# Row 12
$ws.Range("D12").Value = 44553
$ws.Range("M12").Value = 360
$ws.Range("N12").Value = 23000
$ws.Range("O12").Value = 24000
$ws.Range("P12").Value = 23500
$ws.Range("Q12").Value = "`$/caja 16 kilos"
$ws.Range("S12").Value = 1469
$ws.Range("T12").Value = 16

# Row 13
$ws.Range("D13").Value = 44553
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 21000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 21500
$ws.Range("Q13").Value = "`$/caja 16 kilos"
$ws.Range("S13").Value = 1344
$ws.Range("T13").Value = 16

# Row 14
$ws.Range("D14").Value = 44553
$ws.Range("K14").Value = "Modesto"
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 240
$ws.Range("N14").Value = 17000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17500
$ws.Range("Q14").Value = "`$/caja 16 kilos"
$ws.Range("R14").Value = "Región Metropolitana"
$ws.Range("S14").Value = 1094
$ws.Range("T14").Value = 16

# Row 15
$ws.Range("D15").Value = 44573
$ws.Range("K15").Value = "Modesto"
$ws.Range("L15").Value = "Especial"
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 20500
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 20750
$ws.Range("R15").Value = "Región Metropolitana"
$ws.Range("S15").Value = 1153

# Row 16
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 400
$ws.Range("N16").Value = 17500
$ws.Range("O16").Value = 18000
$ws.Range("P16").Value = 17750
$ws.Range("S16").Value = 986

# Row 17
$ws.Range("D17").Value = 44580
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 300
$ws.Range("N17").Value = 22500
$ws.Range("O17").Value = 23000
$ws.Range("P17").Value = 22750
$ws.Range("S17").Value = 1264

# Row 18
$ws.Range("D18").Value = 44580
$ws.Range("K18").Value = "Modesto"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 400
$ws.Range("N18").Value = 19500
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 19750
$ws.Range("R18").Value = "Región Metropolitana"
$ws.Range("S18").Value = 1097

# Row 19
$ws.Range("L19").Value = "Especial"
$ws.Range("M19").Value = 340
$ws.Range("N19").Value = 22500
$ws.Range("O19").Value = 23000
$ws.Range("P19").Value = 22750
$ws.Range("S19").Value = 1264

# Row 20
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 20500
$ws.Range("O20").Value = 21000
$ws.Range("P20").Value = 20750
$ws.Range("S20").Value = 1153

# Row 21
$ws.Range("D21").Value = 44545
$ws.Range("K21").Value = "Castle Brite"
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 300
$ws.Range("N21").Value = 15500
$ws.Range("O21").Value = 16000
$ws.Range("P21").Value = 15750
$ws.Range("Q21").Value = "`$/caja 18 kilos"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 875
$ws.Range("T21").Value = 18

# Row 22
$ws.Range("D22").Value = 44552
$ws.Range("K22").Value = "Castle Brite"
$ws.Range("L22").Value = "Especial"
$ws.Range("M22").Value = 360
$ws.Range("N22").Value = 20000
$ws.Range("O22").Value = 21000
$ws.Range("P22").Value = 20500
$ws.Range("Q22").Value = "`$/caja 18 kilos"
$ws.Range("S22").Value = 1139
$ws.Range("T22").Value = 18

# Row 23
$ws.Range("D23").Value = 44552
$ws.Range("K23").Value = "Castle Brite"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 280
$ws.Range("N23").Value = 18000
$ws.Range("O23").Value = 19000
$ws.Range("P23").Value = 18500
$ws.Range("Q23").Value = "`$/caja 18 kilos"
$ws.Range("S23").Value = 1028
$ws.Range("T23").Value = 18

# Row 24
$ws.Range("D24").Value = 44559
$ws.Range("M24").Value = 400
$ws.Range("N24").Value = 25000
$ws.Range("O24").Value = 26000
$ws.Range("P24").Value = 25500
$ws.Range("S24").Value = 1417

# Row 25
$ws.Range("D25").Value = 44559
$ws.Range("M25").Value = 320
$ws.Range("N25").Value = 22000
$ws.Range("O25").Value = 23000
$ws.Range("P25").Value = 22500
$ws.Range("S25").Value = 1250
